$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "4-" label to the full file name
$ws.Range("B12").Value = "4-write_barrier.txt"

# Fill in the new trial data for the write_barrier row
$ws.Range("D12").Value = "0.101839/0.000000"
$ws.Range("F12").Value = "0.106112/0.000000"
$ws.Range("H12").Value = "0.105315/0.000000"
$ws.Range("C12").Value = "0.891102/0.000000"
$ws.Range("E12").Value = "0.861264/0.000000"
$ws.Range("G12").Value = "0.852867/0.000000"

# Update the current selection to reflect where the user left off
$ws.Range("G12").Select()
